$wb = $excel.ActiveWorkbook

# The "Instructions" sheet holds the step-by-step text in column A.
$ws = $wb.Worksheets.Item("Instructions")

# Row 6 currently reads about saving as a "tab-delimited text file".
# Update it to mention that CSV export is now also supported.
$ws.Range("A6").Value = '4. Save as the excel file (only the "Fill out this form" sheet)  as "Text (tab-delimited) (*.txt)" or "CSV (comma-delimited) (*.csv)"'

# Make the Instructions sheet active and select the edited cell, matching
# the saved workbook view state.
$ws.Activate()
$ws.Range("A6").Select()
